$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full target data for rows 2..17 (descricao, estoque_restante)
# A new row was inserted (CONECTOR FO CAMPO CPO SC APC FRKW) and another
# (DIV_RI 3S 1GHZ H SOLD DESB 488), plus numeric value updates and one
# new row appended at the end (ISOLADOR COAXIAL QUADRADO - CISP-HR),
# shifting subsequent rows down.

$ws.Cells.Item(2, 1).Value = "ACOPLADOR DIREC RI 1S 06DB 1.2G HORIZONT"
$ws.Cells.Item(2, 2).Value = 50

$ws.Cells.Item(3, 1).Value = "CABO FO INVISIVEL SM 1FO G657 B3 RI 50MT"
$ws.Cells.Item(3, 2).Value = 0

$ws.Cells.Item(4, 1).Value = "CONECTOR FO CAMPO CPO SC APC FRKW"
$ws.Cells.Item(4, 2).Value = 82

$ws.Cells.Item(5, 1).Value = "CONTROLE REMOTO BUDGET 00124BU0 00 UEI"
$ws.Cells.Item(5, 2).Value = 69

$ws.Cells.Item(6, 1).Value = "CONTROLE REMOTO VOZ LINUX UEI AVULSO"
$ws.Cells.Item(6, 2).Value = 83

$ws.Cells.Item(7, 1).Value = "CTRL REMOTO DIG. CR2FU UNIVERS"
$ws.Cells.Item(7, 2).Value = 93

$ws.Cells.Item(8, 1).Value = "DIVISOR RI 3S 1.2G HORIZONTAL BALANCEADO"
$ws.Cells.Item(8, 2).Value = 83

$ws.Cells.Item(9, 1).Value = "DIV_RI 3S 1GHZ H SOLD DESB 488"
$ws.Cells.Item(9, 2).Value = 0

$ws.Cells.Item(10, 1).Value = "FITA ACO INOX 1/2 POLEGADA ROLO 25MT"
$ws.Cells.Item(10, 2).Value = 77

$ws.Cells.Item(11, 1).Value = "FITA ACO INOX 3/4 POLEGADA ROLO 25MT"
$ws.Cells.Item(11, 2).Value = 9

$ws.Cells.Item(12, 1).Value = "FONTE ALIM 12V 1.5A ADS18FQ12C12018EPBR"
$ws.Cells.Item(12, 2).Value = 0

$ws.Cells.Item(13, 1).Value = "FONTE ALIM 12V 2.5A MSA C2500IC12030WBR"
$ws.Cells.Item(13, 2).Value = 31

$ws.Cells.Item(14, 1).Value = "FONTE ALIM 12V 3.3A MSA C3330IS12.0 40X"
$ws.Cells.Item(14, 2).Value = 97

$ws.Cells.Item(15, 1).Value = "FONTE ALIM 12V 3A ADS 36FKJ 12 12036EPBR"
$ws.Cells.Item(15, 2).Value = 83

$ws.Cells.Item(16, 1).Value = "FONTE ALIMEN 12V 4A ADS 48PI 12N 212048E"
$ws.Cells.Item(16, 2).Value = 19

$ws.Cells.Item(17, 1).Value = "ISOLADOR COAXIAL QUADRADO - CISP-HR"
$ws.Cells.Item(17, 2).Value = 51
